$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "model3"
$ws.Range("D2").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D6").Value = 1
